# Update "想去人数" (number of people interested) figures on both the
# "展览" sheet and the aggregated "全部类型" sheet.
#
#   F2 : 599  -> 601
#   F8 : 567  -> 574
#   F9 : 3750 -> 3767
#   F10: 68   -> 70

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 601
    $ws.Range("F8").Value = 574
    $ws.Range("F9").Value = 3767
    $ws.Range("F10").Value = 70
}
